{"js": "// Office.js (Word JavaScript API) edit script.\n// This script is the body of `async (context) => { ... }`.\n//\n// Change 1: Collapse the three detailed CORE COMPETENCIES paragraphs into a\n//           single summary paragraph.\n// Change 2: Add a new \"TECHNICAL SKILLS\" section (Heading2 + three body\n//           paragraphs) right after the last \"Trained analytical...\" bullet\n//           and before the closing \"For a more detailed...\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---- Locate the anchor paragraphs by their (stable) text content -------\nlet coreFirstIndex = -1;\nlet coreSecondIndex = -1;\nlet coreThirdIndex = -1;\nlet trainedIndex = -1;\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t.indexOf(\"Product Management & Strategy: Product Conception & Ideation\") === 0) {\n    coreFirstIndex = i;\n  } else if (t.indexOf(\"Technical Product Development: Full-Stack Development\") === 0) {\n    coreSecondIndex = i;\n  } else if (t.indexOf(\"Platform & Infrastructure: Multi-tenant Architecture\") === 0) {\n    coreThirdIndex = i;\n  } else if (t.indexOf(\"Trained analytical and engineering staff on open source geospatial technology\") !== -1) {\n    trainedIndex = i;\n  }\n}\n\nif (coreFirstIndex === -1 || coreSecondIndex === -1 || coreThirdIndex === -1) {\n  throw new Error(\"Could not locate the CORE COMPETENCIES paragraphs.\");\n}\nif (trainedIndex === -1) {\n  throw new Error(\"Could not locate the 'Trained analytical...' paragraph.\");\n}\n\n// ---- Change 1: collapse the 3 CORE COMPETENCIES paragraphs into 1 -------\nitems[coreFirstIndex].insertText(\n  \"Product Management & Strategy \\u2022 Technical Product Development \\u2022 Platform & Infrastructure\",\n  Word.InsertLocation.replace\n);\nitems[coreSecondIndex].delete();\nitems[coreThirdIndex].delete();\nawait context.sync();\n\n// ---- Change 2: insert the new TECHNICAL SKILLS section -------------------\n// Re-fetch the \"Trained analytical...\" paragraph reference fresh (it is\n// still valid since we never touched it), then chain inserts off of it so\n// the three body paragraphs inherit its plain \"Normal\" formatting (no\n// explicit pStyle), and only the heading paragraph gets styled Heading2.\nconst trainedParagraph = items[trainedIndex];\n\nconst platformInfraPara = trainedParagraph.insertParagraph(\n  \"PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Security & Compliance\",\n  Word.InsertLocation.after\n);\nconst techProductDevPara = trainedParagraph.insertParagraph(\n  \"TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; API Development\",\n  Word.InsertLocation.after\n);\nconst productMgmtPara = trainedParagraph.insertParagraph(\n  \"PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development\",\n  Word.InsertLocation.after\n);\nconst headingPara = trainedParagraph.insertParagraph(\"TECHNICAL SKILLS\", Word.InsertLocation.after);\nheadingPara.styleBuiltIn = Word.BuiltInStyleName.heading2;\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change 1: Collapse the three detailed CORE COMPETENCIES paragraphs into a\n#           single summary paragraph.\n# Change 2: Add a new \"TECHNICAL SKILLS\" section (Heading 2 + three body\n#           paragraphs) right after the last \"Trained analytical...\" bullet\n#           and before the closing \"For a more detailed...\" paragraph.\n\n$d = $word.ActiveDocument\n\n# ---- Locate the anchor paragraphs by their (stable) text content --------\n$coreFirstIndex = -1\n$coreSecondIndex = -1\n$coreThirdIndex = -1\n$trainedIndex = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.StartsWith(\"Product Management & Strategy: Product Conception & Ideation\")) {\n        $coreFirstIndex = $i\n    } elseif ($t.StartsWith(\"Technical Product Development: Full-Stack Development\")) {\n        $coreSecondIndex = $i\n    } elseif ($t.StartsWith(\"Platform & Infrastructure: Multi-tenant Architecture\")) {\n        $coreThirdIndex = $i\n    } elseif ($t -like \"*Trained analytical and engineering staff on open source geospatial technology*\") {\n        $trainedIndex = $i\n    }\n}\n\nif ($coreFirstIndex -eq -1 -or $coreSecondIndex -eq -1 -or $coreThirdIndex -eq -1) {\n    throw \"Could not locate the CORE COMPETENCIES paragraphs.\"\n}\nif ($trainedIndex -eq -1) {\n    throw \"Could not locate the 'Trained analytical...' paragraph.\"\n}\n\n# ---- Change 1: collapse the 3 CORE COMPETENCIES paragraphs into 1 -------\n$bullet = [char]0x2022\n$mergedText = \"Product Management & Strategy $bullet Technical Product Development $bullet Platform & Infrastructure\"\n\n$d.Paragraphs.Item($coreFirstIndex).Range.Text = $mergedText\n$d.Paragraphs.Item($coreThirdIndex).Range.Delete()\n$d.Paragraphs.Item($coreSecondIndex).Range.Delete()\n\n# ---- Change 2: insert the new TECHNICAL SKILLS section -------------------\n# The two deletions above shifted every subsequent paragraph's index up by\n# two, so re-locate the \"Trained analytical...\" anchor paragraph fresh\n# instead of reusing the stale $trainedIndex.\n$trainedIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Trained analytical and engineering staff on open source geospatial technology*\") {\n        $trainedIndex = $i\n    }\n}\nif ($trainedIndex -eq -1) {\n    throw \"Could not re-locate the 'Trained analytical...' paragraph after edits.\"\n}\n\n# Chain each insert off the same anchor paragraph (Normal style, no explicit\n# pStyle) so the three body paragraphs inherit plain \"Normal\" formatting and\n# each new paragraph lands immediately after the anchor, in reverse order so\n# the final reading order is correct.\n$anchor = $d.Paragraphs.Item($trainedIndex)\n\n$anchor.Range.InsertParagraphAfter()\n$platformInfraPara = $d.Paragraphs.Item($trainedIndex + 1)\n$platformInfraPara.Range.Text = \"PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Security & Compliance\"\n\n$anchor.Range.InsertParagraphAfter()\n$techProductDevPara = $d.Paragraphs.Item($trainedIndex + 1)\n$techProductDevPara.Range.Text = \"TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; API Development\"\n\n$anchor.Range.InsertParagraphAfter()\n$productMgmtPara = $d.Paragraphs.Item($trainedIndex + 1)\n$productMgmtPara.Range.Text = \"PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development\"\n\n$anchor.Range.InsertParagraphAfter()\n$headingPara = $d.Paragraphs.Item($trainedIndex + 1)\n$headingPara.Range.Text = \"TECHNICAL SKILLS\"\n$headingPara.Style = \"Heading 2\"\n"}
